$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160:237 down to 161:238
$ws.Rows.Item(160).Insert()

# Populate the newly-inserted row 160 with the new record's data.
$ws.Cells.Item(160, 1).Value = 3
$ws.Cells.Item(160, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(160, 3).Value = "Coquimbo"
$ws.Cells.Item(160, 4).Value = 44466
$ws.Cells.Item(160, 5).Value = 5
$ws.Cells.Item(160, 6).Value = 100112017
$ws.Cells.Item(160, 7).Value = "Apio"
$ws.Cells.Item(160, 8).Value = "Americana (o)"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 230
$ws.Cells.Item(160, 11).Value = 9000
$ws.Cells.Item(160, 12).Value = 9500
$ws.Cells.Item(160, 13).Value = 9239
$ws.Cells.Item(160, 14).Value = "$/docena de matas"
$ws.Cells.Item(160, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(160, 16).Value = 1540
$ws.Cells.Item(160, 17).Value = 6
$ws.Cells.Item(160, 18).Value = "Hortaliza"
